$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Marzo de 2020 a las 16:46"

$ws.Cells.Item(7, 2).Value = 27151
$ws.Cells.Item(7, 5).Value = 26624
$ws.Cells.Item(7, 3).Value = 2944

$ws.Cells.Item(8, 2).Value = 23974
$ws.Cells.Item(8, 5).Value = 23615
$ws.Cells.Item(8, 3).Value = 1610

$ws.Cells.Item(13, 7).Value = 11
$ws.Cells.Item(13, 8).Value = 244
$ws.Cells.Item(13, 5).Value = 4681

$ws.Cells.Item(16, 7).Value = 8
$ws.Cells.Item(16, 5).Value = 3277
$ws.Cells.Item(16, 8).Value = 16

$ws.Cells.Item(17, 5).Value = 2238
$ws.Cells.Item(17, 3).Value = 87
$ws.Cells.Item(17, 2).Value = 2251

$ws.Cells.Item(25, 5).Value = 815
$ws.Cells.Item(25, 3).Value = 32
$ws.Cells.Item(25, 2).Value = 1086

$ws.Cells.Item(35, 1).Value = "Grecia"
$ws.Cells.Item(35, 6).Value = 18
$ws.Cells.Item(35, 4).Value = 19
$ws.Cells.Item(35, 8).Value = 13
$ws.Cells.Item(35, 5).Value = 592
$ws.Cells.Item(35, 3).Value = 94
$ws.Cells.Item(35, 2).Value = 624

$ws.Cells.Item(36, 1).Value = "Tailandia"
$ws.Cells.Item(36, 4).Value = 44
$ws.Cells.Item(36, 6).Value = 7
$ws.Cells.Item(36, 5).Value = 554
$ws.Cells.Item(36, 2).Value = 599
$ws.Cells.Item(36, 3).Value = 188

$ws.Cells.Item(37, 1).Value = "Islandia"
$ws.Cells.Item(37, 4).Value = 5
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 2).Value = 568
$ws.Cells.Item(37, 8).Value = 1
$ws.Cells.Item(37, 3).Value = 95
$ws.Cells.Item(37, 5).Value = 562
$ws.Cells.Item(37, 6).Value = 1

$ws.Cells.Item(38, 1).Value = "Polonia"
$ws.Cells.Item(38, 6).Value = 3
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 3).Value = 27
$ws.Cells.Item(38, 4).Value = 13
$ws.Cells.Item(38, 5).Value = 543
$ws.Cells.Item(38, 2).Value = 563

$ws.Cells.Item(39, 1).Value = "Ecuador"
$ws.Cells.Item(39, 5).Value = 522
$ws.Cells.Item(39, 6).Value = 2
$ws.Cells.Item(39, 4).Value = 3
$ws.Cells.Item(39, 2).Value = 532
$ws.Cells.Item(39, 8).Value = 7

$ws.Cells.Item(54, 1).Value = "Croacia"
$ws.Cells.Item(54, 6).Value = 5
$ws.Cells.Item(54, 2).Value = 254
$ws.Cells.Item(54, 5).Value = 248
$ws.Cells.Item(54, 8).Value = 1
$ws.Cells.Item(54, 4).Value = 5

$ws.Cells.Item(55, 1).Value = "Mexico"
$ws.Cells.Item(55, 3).Value = 48
$ws.Cells.Item(55, 6).Value = 1
$ws.Cells.Item(55, 8).Value = 2
$ws.Cells.Item(55, 2).Value = 251
$ws.Cells.Item(55, 5).Value = 245
$ws.Cells.Item(55, 4).Value = 4

$ws.Cells.Item(56, 1).Value = "Libano"
$ws.Cells.Item(56, 2).Value = 248
$ws.Cells.Item(56, 8).Value = 4
$ws.Cells.Item(56, 4).Value = 8
$ws.Cells.Item(56, 5).Value = 236
$ws.Cells.Item(56, 3).Value = 18
$ws.Cells.Item(56, 6).Value = 4

$ws.Cells.Item(57, 1).Value = "Panama"
$ws.Cells.Item(57, 5).Value = 241
$ws.Cells.Item(57, 2).Value = 245
$ws.Cells.Item(57, 6).Value = 7
$ws.Cells.Item(57, 4).Value = 1
$ws.Cells.Item(57, 8).Value = 3

$ws.Cells.Item(58, 1).Value = "Sudafrica"
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 2).Value = 240
$ws.Cells.Item(58, 5).Value = 238
$ws.Cells.Item(58, 4).Value = 2
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 3).Value = 0

$ws.Cells.Item(60, 1).Value = "Colombia"
$ws.Cells.Item(60, 2).Value = 231
$ws.Cells.Item(60, 8).Value = 2
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 3).Value = 35
$ws.Cells.Item(60, 5).Value = 226
$ws.Cells.Item(60, 4).Value = 3

$ws.Cells.Item(61, 1).Value = "Argentina"
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 4).Value = 27
$ws.Cells.Item(61, 3).Value = 67
$ws.Cells.Item(61, 2).Value = 225
$ws.Cells.Item(61, 8).Value = 4
$ws.Cells.Item(61, 5).Value = 194

$ws.Cells.Item(74, 1).Value = "Hungria"
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 4).Value = 16
$ws.Cells.Item(74, 8).Value = 6
$ws.Cells.Item(74, 6).Value = 6
$ws.Cells.Item(74, 5).Value = 109
$ws.Cells.Item(74, 3).Value = 28

$ws.Cells.Item(75, 1).Value = "Lituania"
$ws.Cells.Item(75, 8).Value = 1
$ws.Cells.Item(75, 5).Value = 127
$ws.Cells.Item(75, 2).Value = 129
$ws.Cells.Item(75, 4).Value = 1
$ws.Cells.Item(75, 6).Value = 1
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 3).Value = 30

$ws.Cells.Item(83, 1).Value = "Moldavia"
$ws.Cells.Item(83, 3).Value = 14
$ws.Cells.Item(83, 6).Value = 3
$ws.Cells.Item(83, 5).Value = 92
$ws.Cells.Item(83, 4).Value = 1

$ws.Cells.Item(84, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(84, 8).Value = 1
$ws.Cells.Item(84, 5).Value = 91
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(84, 2).Value = 94

$ws.Cells.Item(85, 1).Value = "Malta"
$ws.Cells.Item(85, 2).Value = 90
$ws.Cells.Item(85, 3).Value = 17
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 6).Value = 1
$ws.Cells.Item(85, 5).Value = 88

$ws.Cells.Item(86, 1).Value = "Albania"
$ws.Cells.Item(86, 3).Value = 13
$ws.Cells.Item(86, 2).Value = 89
$ws.Cells.Item(86, 8).Value = 2
$ws.Cells.Item(86, 5).Value = 85

$ws.Cells.Item(87, 1).Value = "Brunei"
$ws.Cells.Item(87, 6).Value = 2
$ws.Cells.Item(87, 4).Value = 2
$ws.Cells.Item(87, 2).Value = 88
$ws.Cells.Item(87, 5).Value = 86
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 3).Value = 5

$ws.Cells.Item(88, 1).Value = "Camboya"
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 4).Value = 2
$ws.Cells.Item(88, 5).Value = 82
$ws.Cells.Item(88, 2).Value = 84
$ws.Cells.Item(88, 3).Value = 31

$ws.Cells.Item(89, 1).Value = "Republica de Chipre"
$ws.Cells.Item(89, 4).Value = 3
$ws.Cells.Item(89, 2).Value = 84
$ws.Cells.Item(89, 5).Value = 80

$ws.Cells.Item(90, 1).Value = "Sri Lanka"
$ws.Cells.Item(90, 3).Value = 5
$ws.Cells.Item(90, 2).Value = 82
$ws.Cells.Item(90, 6).Value = 2
$ws.Cells.Item(90, 4).Value = 3
$ws.Cells.Item(90, 5).Value = 79

$ws.Cells.Item(91, 1).Value = "Bielorrusia"
$ws.Cells.Item(91, 2).Value = 76
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 5).Value = 61
$ws.Cells.Item(91, 4).Value = 15
$ws.Cells.Item(91, 7).Value = 0

$ws.Cells.Item(92, 1).Value = "Tunez"
$ws.Cells.Item(92, 3).Value = 15
$ws.Cells.Item(92, 8).Value = 3
$ws.Cells.Item(92, 7).Value = 2
$ws.Cells.Item(92, 5).Value = 71
$ws.Cells.Item(92, 6).Value = 7
$ws.Cells.Item(92, 4).Value = 1

$ws.Cells.Item(93, 1).Value = "Burkina Faso"
$ws.Cells.Item(93, 5).Value = 66
$ws.Cells.Item(93, 2).Value = 75
$ws.Cells.Item(93, 8).Value = 4
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 3).Value = 11
$ws.Cells.Item(93, 4).Value = 5
$ws.Cells.Item(93, 6).Value = 0

$ws.Cells.Item(94, 1).Value = "Venezuela"
$ws.Cells.Item(94, 4).Value = 15
$ws.Cells.Item(94, 2).Value = 70
$ws.Cells.Item(94, 3).Value = 0
$ws.Cells.Item(94, 6).Value = 2
$ws.Cells.Item(94, 5).Value = 55

$ws.Cells.Item(95, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(95, 5).Value = 66
$ws.Cells.Item(95, 3).Value = 14
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 2).Value = 66
$ws.Cells.Item(95, 4).Value = 0

$ws.Cells.Item(96, 1).Value = "Azerbaiyan"
$ws.Cells.Item(96, 8).Value = 1
$ws.Cells.Item(96, 5).Value = 53
$ws.Cells.Item(96, 2).Value = 65
$ws.Cells.Item(96, 3).Value = 12
$ws.Cells.Item(96, 4).Value = 11

$ws.Cells.Item(97, 3).Value = 5
$ws.Cells.Item(97, 5).Value = 59
$ws.Cells.Item(97, 2).Value = 59

$ws.Cells.Item(98, 1).Value = "Estado de Palestina"
$ws.Cells.Item(98, 3).Value = 6
$ws.Cells.Item(98, 4).Value = 17
$ws.Cells.Item(98, 5).Value = 42
$ws.Cells.Item(98, 2).Value = 59
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 8).Value = 0

$ws.Cells.Item(99, 1).Value = "Guadalupe"
$ws.Cells.Item(99, 6).Value = 4
$ws.Cells.Item(99, 5).Value = 55
$ws.Cells.Item(99, 8).Value = 1
$ws.Cells.Item(99, 4).Value = 0

$ws.Cells.Item(100, 1).Value = "Senegal"
$ws.Cells.Item(100, 4).Value = 5
$ws.Cells.Item(100, 2).Value = 56
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 5).Value = 51

$ws.Cells.Item(101, 1).Value = "Oman"
$ws.Cells.Item(101, 2).Value = 55
$ws.Cells.Item(101, 3).Value = 3
$ws.Cells.Item(101, 5).Value = 38
$ws.Cells.Item(101, 4).Value = 17
$ws.Cells.Item(101, 6).Value = 0

$ws.Cells.Item(102, 1).Value = "Georgia"
$ws.Cells.Item(102, 4).Value = 1
$ws.Cells.Item(102, 2).Value = 54
$ws.Cells.Item(102, 5).Value = 53
$ws.Cells.Item(102, 6).Value = 1
$ws.Cells.Item(102, 3).Value = 5

$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 5).Value = 32
$ws.Cells.Item(110, 8).Value = 1

$ws.Cells.Item(177, 1).Value = "Sudan"
$ws.Cells.Item(177, 8).Value = 1
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 3).Value = 0

$ws.Cells.Item(178, 1).Value = "Montserrat"
$ws.Cells.Item(178, 2).Value = 1
$ws.Cells.Item(178, 8).Value = 0

$ws.Cells.Item(179, 1).Value = "Gambia"

$ws.Cells.Item(180, 1).Value = "San Vicente y las Granadinas"

$ws.Cells.Item(181, 1).Value = "Santa Sede"

$ws.Cells.Item(182, 1).Value = "Eritrea"

$ws.Cells.Item(183, 1).Value = "Papua Nueva Guinea"

$ws.Cells.Item(184, 1).Value = "Uganda"

$ws.Cells.Item(185, 1).Value = "Somalia"

$ws.Cells.Item(186, 1).Value = "Republica de Yibuti"

$ws.Cells.Item(187, 1).Value = "Timor Oriental"

$ws.Cells.Item(188, 1).Value = "Republica del Chad"

$ws.Cells.Item(189, 1).Value = "Antigua y Barbuda"

$ws.Cells.Item(190, 1).Value = "San Martin (Parte Holandesa)"

$ws.Cells.Item(191, 1).Value = "Nepal"
$ws.Cells.Item(191, 4).Value = 1
$ws.Cells.Item(191, 5).Value = 0
